$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "log" entry (row 13): copy formatting down from the row above (row 12)
# then overwrite the values, matching what Excel does when a user continues
# typing new rows below an existing table.
$ws.Range("D12").Copy($ws.Range("D13"))
$ws.Range("D13").Value = 44260

$ws.Range("E12").Copy($ws.Range("E13"))
$ws.Range("E13").Value = "Résolution d'un bug pour les grilles"

$ws.Rows.Item(13).RowHeight = 30

[void]$ws.Range("E14").Select()
